# Set column C ("showPercentage") to 0 for all variable rows that don't
# already have a value there (rows 7 and 8 already contain 1 and are left
# untouched). This mirrors the diff which adds <c r="C2">0</c> ... <c r="C20">0</c>
# (skipping C7/C8) to xl/worksheets/sheet1.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,9,10,11,12,13,14,15,16,17,18,19,20)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 0
}
